$d = $word.ActiveDocument

# Locate the run of text "PAQJP_3" that needs to become two runs:
# "PAQJP_" and "4" (both keeping the same "s2" character style).
$rng = $d.Content
$found = $rng.Find.Execute("PAQJP_3", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    # $rng now spans the matched "PAQJP_3" text (Start..End).
    # Remove the trailing "3" character first.
    $threeRange = $d.Range($rng.End - 1, $rng.End)
    $threeRange.Delete()

    # Insert a new "4" right after "PAQJP_" as its own run, then give
    # that new run the same "s2" character style as the rest of the text.
    $insertPos = $rng.End - 1
    $newRng = $d.Range($insertPos, $insertPos)
    $newRng.InsertAfter("4")
    $newRng.Style = "s2"
}
